# Update statistical values on the "Tab31" worksheet (informal employment /
# informal economy indicators). The underlying source data for Gambia (GMB,
# row 50) was revised, which in turn changes every aggregate/group row that
# includes Gambia (Afrique de l'Ouest, Afrique, Afrique hors ressources,
# PMA, Etats fragiles, etc.) for the first three indicator columns (C, D, E).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab31")

# --- Gambie (GMB) ---------------------------------------------------------
$ws.Range("C50").Value = 81.599999999999994
$ws.Range("D50").Value = 86.9
$ws.Range("E50").Value = 76

# --- Afrique de l'Ouest -----------------------------------------------------
$ws.Range("C61").Value = 87.65
$ws.Range("D61").Value = 90.985714285714295
$ws.Range("E61").Value = 84.8642857142857

# --- Afrique -----------------------------------------------------------------
$ws.Range("C62").Value = 81.897560975609807
$ws.Range("D62").Value = 83.921951219512195
$ws.Range("E62").Value = 79.914634146341498

# --- Reste du monde ----------------------------------------------------------
$ws.Range("C63").Value = 40.762790697674397
$ws.Range("D63").Value = 40.327906976744202
$ws.Range("E63").Value = 41.0162790697674

# --- Amerique latine et Caraibes ---------------------------------------------
$ws.Range("C64").Value = 56.359090909090902
$ws.Range("D64").Value = 54.731818181818198
$ws.Range("E64").Value = 57.572727272727299

# --- Monde ---------------------------------------------------------------------
$ws.Range("C66").Value = 54.042519685039402
$ws.Range("D66").Value = 54.4015748031496
$ws.Range("E66").Value = 53.574015748031499

# --- CEDEAO ----------------------------------------------------------------------
$ws.Range("C68").Value = 83.94
$ws.Range("D68").Value = 86.1
$ws.Range("E68").Value = 81.915000000000006

# --- Afrique, pays riches en ressources -----------------------------------------
$ws.Range("C71").Value = 87.65
$ws.Range("D71").Value = 90.985714285714295
$ws.Range("E71").Value = 84.8642857142857

# --- Afrique (pays riches en ressources exclus) ---------------------------------
$ws.Range("C82").Value = 81.5513513513514
$ws.Range("D82").Value = 83.5324324324325
$ws.Range("E82").Value = 79.583783783783801

# --- RDM (pays riches en ressources exclus) -------------------------------------
$ws.Range("C83").Value = 38.792307692307702
$ws.Range("D83").Value = 38.729487179487201
$ws.Range("E83").Value = 38.9

# --- Afrique, pays a faible revenu ----------------------------------------------
$ws.Range("C84").Value = 91.0833333333333
$ws.Range("D84").Value = 93.938888888888897
$ws.Range("E84").Value = 88.605555555555597

# --- Afrique, pays les moins avances --------------------------------------------
$ws.Range("C90").Value = 14.0432432432432
$ws.Range("D90").Value = 13.4
$ws.Range("E90").Value = 14.524324324324301

# --- RDM, pays les moins avances -------------------------------------------------
$ws.Range("C91").Value = 89.040740740740802
$ws.Range("D91").Value = 92.248148148148204
$ws.Range("E91").Value = 86.551851851851893

# --- Afrique, Etats fragiles -----------------------------------------------------
$ws.Range("C97").Value = 86.969230769230805
$ws.Range("D97").Value = 90.838461538461601
$ws.Range("E97").Value = 83.873076923076894
